# 物品表增加id — add a new "ItemId" (int) column to the "resource" sheet.
#
# The "resource" sheet is a config-table style sheet laid out as:
#   row1 = ##var   (english field names)
#   row2 = ##type  (field types)
#   row3 = ##group (grouping marker)
#   row4 = ##      (chinese field names)
#   row5+ = actual data rows (one per resource: Coin, Crystal, Flower, Leaf, Money, Other)
#
# A new column F is introduced for the item id:
#   F1 = ItemId   (var name)
#   F2 = int      (type)
#   F4 = 物品编号  (chinese name)
#   F5:F10 = 1..6 (one id per resource row)

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("resource")

$ws1.Range("F1").Value = "ItemId"
$ws1.Range("F2").Value = "int"
$ws1.Range("F4").Value = "物品编号"

$ws1.Range("F5").Value  = 1
$ws1.Range("F6").Value  = 2
$ws1.Range("F7").Value  = 3
$ws1.Range("F8").Value  = 4
$ws1.Range("F9").Value  = 5
$ws1.Range("F10").Value = 6

# Match the author's final cursor position on the "resource" sheet.
[void]$ws1.Range("D7").Select()
